# Mock trial data expanded: team roster entry that used to read "Team61"
# is replaced with "Team16" (16th team added as the mock trial grew).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Team16"

# Reset the worksheet's active cell/selection back to the top-left corner
# (A1) now that the data entry checks are done, instead of leaving the
# selection parked wherever it last was (I12).
$ws.Range("A1").Select()
